$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.769|x12:0.000|x13:0.231|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B3").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.180|x9:0.621|x10:0.199|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B4").Value = "x1:0.000|x2:0.186|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.628|x13:0.000|x14:0.000|x15:0.000|x16:0.186|x17:0.000"
$ws.Range("B5").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:1.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B7").Value = "x1:0.000|x2:0.495|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.505|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B10").Value = "x1:0.509|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.295|x9:0.000|x10:0.196|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B11").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.339|x12:0.661|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B12").Value = "x1:0.393|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.607|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B13").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.439|x8:0.000|x9:0.220|x10:0.000|x11:0.000|x12:0.000|x13:0.341|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B16").Value = "x1:0.250|x2:0.000|x3:0.000|x4:0.750|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B18").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.273|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.727|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B19").Value = "x1:1.152|x2:0.681|x3:1.000|x4:0.750|x5:1.000|x6:1.000|x7:1.216|x8:1.475|x9:0.842|x10:1.001|x11:1.107|x12:1.289|x13:1.300|x14:1.000|x15:1.000|x16:1.186|x17:0.000"
